$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9750066666666667
$ws.Range("H2").Value = 2.92502
$ws.Range("I2").Value = 0.03314938328005726
$ws.Range("J2").Value = 0.03314938328005725
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.217066
$ws.Range("N2").Value = 0.6511979999999999
$ws.Range("O2").Value = 0.06317054032028298
$ws.Range("P2").Value = 0.06317054032028299
$ws.Range("Q2").Value = 0.2116407971066666
$ws.Range("R2").Value = 1.90476717396
$ws.Range("S2").Value = 0.002094064453085371
$ws.Range("T2").Value = 0.002094064453085371
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9750066666666667
$ws.Range("H3").Value = 2.92502
$ws.Range("I3").Value = 0.03314938328005726
$ws.Range("J3").Value = 0.03314938328005725
$ws.Range("O3").Value = 0.1114685051885465
$ws.Range("P3").Value = 0.1114685051885465
$ws.Range("Q3").Value = 0.3734538785133333
$ws.Range("R3").Value = 3.36108490662
$ws.Range("S3").Value = 0.003695112202150179
$ws.Range("T3").Value = 0.003695112202150178
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9750066666666667
$ws.Range("H4").Value = 2.92502
$ws.Range("I4").Value = 0.03314938328005726
$ws.Range("J4").Value = 0.03314938328005725
$ws.Range("M4").Value = 1.344541
$ws.Range("N4").Value = 4.033623
$ws.Range("O4").Value = 0.3912882784626502
$ws.Range("P4").Value = 0.3912882784626502
$ws.Range("Q4").Value = 1.310936438606667
$ws.Range("R4").Value = 11.79842794746
$ws.Range("S4").Value = 0.01297096511575217
$ws.Range("T4").Value = 0.01297096511575216
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9750066666666667
$ws.Range("H5").Value = 2.92502
$ws.Range("I5").Value = 0.03314938328005726
$ws.Range("J5").Value = 0.03314938328005725
$ws.Range("M5").Value = 0.3057526666666667
$ws.Range("N5").Value = 0.917258
$ws.Range("O5").Value = 0.08898013119374158
$ws.Range("P5").Value = 0.0889801311937416
$ws.Range("Q5").Value = 0.2981108883511111
$ws.Range("R5").Value = 2.68299799516
$ws.Range("S5").Value = 0.002949636473251119
$ws.Range("T5").Value = 0.002949636473251118
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9750066666666667
$ws.Range("H6").Value = 2.92502
$ws.Range("I6").Value = 0.03314938328005726
$ws.Range("J6").Value = 0.03314938328005725
$ws.Range("M6").Value = 1.185803666666667
$ws.Range("N6").Value = 3.557411
$ws.Range("O6").Value = 0.3450925448347787
$ws.Range("P6").Value = 0.3450925448347787
$ws.Range("Q6").Value = 1.156166480357778
$ws.Range("R6").Value = 10.40549832322
$ws.Range("S6").Value = 0.01143960503581842
$ws.Range("T6").Value = 0.01143960503581842
$ws.Range("I7").Value = 0.9153383232957121
$ws.Range("J7").Value = 0.9153383232957121
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.217066
$ws.Range("N7").Value = 0.6511979999999999
$ws.Range("O7").Value = 0.06317054032028298
$ws.Range("P7").Value = 0.06317054032028299
$ws.Range("Q7").Value = 5.843937750755333
$ws.Range("R7").Value = 52.59543975679799
$ws.Range("S7").Value = 0.05782241645845199
$ws.Range("T7").Value = 0.05782241645845201
$ws.Range("I8").Value = 0.9153383232957121
$ws.Range("J8").Value = 0.9153383232957121
$ws.Range("O8").Value = 0.1114685051885465
$ws.Range("P8").Value = 0.1114685051885465
$ws.Range("S8").Value = 0.1020313946395635
$ws.Range("T8").Value = 0.1020313946395635
$ws.Range("I9").Value = 0.9153383232957121
$ws.Range("J9").Value = 0.9153383232957121
$ws.Range("M9").Value = 1.344541
$ws.Range("N9").Value = 4.033623
$ws.Range("O9").Value = 0.3912882784626502
$ws.Range("P9").Value = 0.3912882784626502
$ws.Range("Q9").Value = 36.19827106658034
$ws.Range("R9").Value = 325.784439599223
$ws.Range("S9").Value = 0.3581611567332679
$ws.Range("T9").Value = 0.3581611567332679
$ws.Range("I10").Value = 0.9153383232957121
$ws.Range("J10").Value = 0.9153383232957121
$ws.Range("M10").Value = 0.3057526666666667
$ws.Range("N10").Value = 0.917258
$ws.Range("O10").Value = 0.08898013119374158
$ws.Range("P10").Value = 0.0889801311937416
$ws.Range("Q10").Value = 8.231595694984222
$ws.Range("R10").Value = 74.084361254858
$ws.Range("S10").Value = 0.08144692409351191
$ws.Range("T10").Value = 0.08144692409351192
$ws.Range("I11").Value = 0.9153383232957121
$ws.Range("J11").Value = 0.9153383232957121
$ws.Range("M11").Value = 1.185803666666667
$ws.Range("N11").Value = 3.557411
$ws.Range("O11").Value = 0.3450925448347787
$ws.Range("P11").Value = 0.3450925448347787
$ws.Range("Q11").Value = 31.92468103073455
$ws.Range("R11").Value = 287.322129276611
$ws.Range("S11").Value = 0.3158764313709167
$ws.Range("T11").Value = 0.3158764313709167
$ws.Range("G12").Value = 1.515106
$ws.Range("H12").Value = 4.545318
$ws.Range("I12").Value = 0.05151229342423071
$ws.Range("J12").Value = 0.0515122934242307
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.217066
$ws.Range("N12").Value = 0.6511979999999999
$ws.Range("O12").Value = 0.06317054032028298
$ws.Range("P12").Value = 0.06317054032028299
$ws.Range("Q12").Value = 0.328877998996
$ws.Range("R12").Value = 2.959901990964
$ws.Range("S12").Value = 0.003254059408745613
$ws.Range("T12").Value = 0.003254059408745614
$ws.Range("G13").Value = 1.515106
$ws.Range("H13").Value = 4.545318
$ws.Range("I13").Value = 0.05151229342423071
$ws.Range("J13").Value = 0.0515122934242307
$ws.Range("O13").Value = 0.1114685051885465
$ws.Range("P13").Value = 0.1114685051885465
$ws.Range("Q13").Value = 0.5803265058620001
$ws.Range("R13").Value = 5.222938552758
$ws.Range("S13").Value = 0.00574199834683279
$ws.Range("T13").Value = 0.005741998346832789
$ws.Range("G14").Value = 1.515106
$ws.Range("H14").Value = 4.545318
$ws.Range("I14").Value = 0.05151229342423071
$ws.Range("J14").Value = 0.0515122934242307
$ws.Range("M14").Value = 1.344541
$ws.Range("N14").Value = 4.033623
$ws.Range("O14").Value = 0.3912882784626502
$ws.Range("P14").Value = 0.3912882784626502
$ws.Range("Q14").Value = 2.037122136346
$ws.Range("R14").Value = 18.334099227114
$ws.Range("S14").Value = 0.02015615661363013
$ws.Range("T14").Value = 0.02015615661363013
$ws.Range("G15").Value = 1.515106
$ws.Range("H15").Value = 4.545318
$ws.Range("I15").Value = 0.05151229342423071
$ws.Range("J15").Value = 0.0515122934242307
$ws.Range("M15").Value = 0.3057526666666667
$ws.Range("N15").Value = 0.917258
$ws.Range("O15").Value = 0.08898013119374158
$ws.Range("P15").Value = 0.0889801311937416
$ws.Range("Q15").Value = 0.4632476997826667
$ws.Range("R15").Value = 4.169229298044
$ws.Range("S15").Value = 0.00458357062697856
$ws.Range("T15").Value = 0.004583570626978561
$ws.Range("G16").Value = 1.515106
$ws.Range("H16").Value = 4.545318
$ws.Range("I16").Value = 0.05151229342423071
$ws.Range("J16").Value = 0.0515122934242307
$ws.Range("M16").Value = 1.185803666666667
$ws.Range("N16").Value = 3.557411
$ws.Range("O16").Value = 0.3450925448347787
$ws.Range("P16").Value = 0.3450925448347787
$ws.Range("Q16").Value = 1.796618250188667
$ws.Range("R16").Value = 10.40549832322
$ws.Range("S16").Value = 0.01777650842804361
$ws.Range("T16").Value = 0.01777650842804361
